# Daily scrape update: refresh Global Talent opportunity listing
# - replaces the 3 existing data rows (2-4) with new scraped data
# - appends 11 new rows (5-15)
# - highlights the PREMIUM column with a yellow fill when the value is "Yes"
# - widens several columns to fit the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row data: Id, Title, Country, Premium, Applicants, Duration, Organization
$rows = @(
  @("1327884", "[EXP] Customer Emission Reporting (EU Preferred)", "Fritz-Erler-Straße 5, 53113 Bonn, Germany", "Yes", "5 applicants", "6 - 18 Months", "DHL Group"),
  @("1327883", "Mobile Application Developer", "El-Mahalla El-Kubra, Al Mahalah Al Kubra (Part 2), El Mahalla El Kubra, Gharbia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Positive Kids academy"),
  @("1327882", "Web developer", "El-Mahalla El-Kubra, Al Mahalah Al Kubra (Part 2), El Mahalla El Kubra, Gharbia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Positive Kids academy"),
  @("1327880", "Legal Support - Procurement Sustainability Department (EU Only)", "Madrid, Spain", "Yes", "1 applicant", "6 - 18 Months", "DHL Group"),
  @("1327878", "Legal Support - Procurement Sustainability Department (EU Only)", "Bruxelles, Belgium", "Yes", "4 applicants", "6 - 18 Months", "DHL Group"),
  @("1327863", "Customer Support", "8200 Albufeira, Portugal", "No", "9 applicants", "9 - 12 Weeks", "BERNARDINO GOMES - GESTÃO HOTELEIRA, S.A."),
  @("1327810", "Photographer", "El Sadat City, Menofia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Habib Agency"),
  @("1327809", "Video editor", "El Sadat City, Menofia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Habib Agency"),
  @("1327808", "Graphic designer", "El Sadat City, Menofia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Habib Agency"),
  @("1327805", "Content creator", "Shebeen El-Kom, Qism Shebeen El-Kom, Shibin el Kom, Menofia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "European Hospital"),
  @("1327644", "Marketing Executive Intern", "Chandigarh, India", "No", "1 applicant", "9 - 12 Weeks", "Silver Fern Education Consultants"),
  @("1324106", "Tourism Specialist - Intern", "Nugegoda, Sri Lanka", "No", "18 applicants", "9 - 12 Weeks", "Brand Corridor (Pvt) Ltd"),
  @("1322493", "[Impact Fortaleza]- Cost & Quality Planning", "Castanhal, PA, Brasil", "No", "26 applicants", "6 - 18 Months", "Petruz Fruity"),
  @("1321451", "Sales & Marketing Representative", "Edmonton, AB, Canada", "No", "51 applicants", "6 - 18 Months", "Canada Prime Marketing")
)

$r = 2
foreach ($row in $rows) {
    $id = $row[0]
    $url = "https://aiesec.org/opportunity/global-talent/" + $id

    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $url
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]

    if ($row[3] -eq "Yes") {
        $ws.Cells.Item($r, 5).Interior.Color = 65535
    }

    $r++
}

# Column width adjustments to fit the refreshed content
$ws.Columns.Item(3).ColumnWidth = 65.16666666666667
$ws.Columns.Item(4).ColumnWidth = 100.16666666666667
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 43.166666666666664
